# Add the flood warning map, issue #168
#
# This updates the flood-orgs.xlsx workbook:
#  - Rewrites the "Notes" sheet description text.
#  - Reorders/rewrites the "Data" sheet header row and existing data rows
#    (columns were re-sequenced: Organization, OrganizationType, Description,
#    FloodWarningWebpage, FloodplainRegWebpage, InBasin, Longitude, Latitude,
#    LongitudeMain, LatitudeMain, Note).
#  - Adds a new hyperlink on D5 (Larimer County FloodWarningWebpage).
#  - Adds a new organization row (Boxelder Regional Stormwater Authority).
#  - Updates saved sheet view selections.

$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")
$data  = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# Data sheet: header row (columns re-sequenced).
# ---------------------------------------------------------------------
$data.Range("A1").Value = "Organization"
$data.Range("B1").Value = "OrganizationType"
$data.Range("C1").Value = "Description"
$data.Range("D1").Value = "FloodWarningWebpage"
$data.Range("E1").Value = "FloodplainRegWebpage"
$data.Range("F1").Value = "InBasin"
$data.Range("G1").Value = "Longitude"
$data.Range("H1").Value = "Latitude"
$data.Range("I1").Value = "LongitudeMain"
$data.Range("J1").Value = "LatitudeMain"
$data.Range("K1").Value = "Note"

# ---------------------------------------------------------------------
# Data sheet: row 2 - City of Greeley
# ---------------------------------------------------------------------
$data.Range("A2").Value = "City of Greeley"
$data.Range("B2").Value = "Municipal Utility"
$data.Range("C2").Value = "Stormwater - Floodplain Information"
$data.Range("D2").Value = "https://larimerco-ns5.trilynx-novastar.systems/novastar/operator"
$data.Range("E2").Value = "https://greeleygov.com/services/pw/stormwater/floodplain-information"
$data.Range("F2").Value = "Yes"

# ---------------------------------------------------------------------
# Data sheet: row 3 - City of Fort Collins
# ---------------------------------------------------------------------
$data.Range("A3").Value = "City of Fort Collins"
$data.Range("B3").Value = "Municipal Utility"
$data.Range("C3").Value = "Utilities - Flooding"
$data.Range("D3").Value = "https://www.fcgov.com/utilities/what-we-do/stormwater/flooding/warning-system"
$data.Range("E3").Value = "https://www.fcgov.com/utilities/what-we-do/stormwater/flooding/"
$data.Range("F3").Value = "Yes"

# ---------------------------------------------------------------------
# Data sheet: row 4 - City of Loveland
# ---------------------------------------------------------------------
$data.Range("A4").Value = "City of Loveland"
$data.Range("B4").Value = "Municipal Utility"
$data.Range("C4").Value = "Stormwater - Flood Management"
$data.Range("D4").Value = "https://larimerco-ns5.trilynx-novastar.systems/novastar/operator"
$data.Range("E4").Value = "https://www.lovgov.org/services/public-works/stormwater/flood-management"
$data.Range("F4").Value = "Yes"

# ---------------------------------------------------------------------
# Data sheet: row 5 - Larimer County (gains a FloodWarningWebpage hyperlink)
# ---------------------------------------------------------------------
$data.Range("A5").Value = "Larimer County"
$data.Range("B5").Value = "County Government"
$data.Range("C5").Value = "Floodplains"
$data.Range("D5").Value = "https://larimerco-ns5.trilynx-novastar.systems/novastar/operator"
$data.Range("E5").Value = "https://www.larimer.org/engineering/floodplains"
$data.Range("F5").Value = "Yes"
$data.Hyperlinks.Add($data.Range("D5"), "https://larimerco-ns5.trilynx-novastar.systems/novastar/operator")
$data.Range("D5").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Data sheet: row 6 - new organization, Boxelder Regional Stormwater Authority
# ---------------------------------------------------------------------
$data.Range("A6").Value = "Boxelder Regional Stormwater Authority"
$data.Range("B6").Value = "IGA"
$data.Range("C6").Value = "Stormwater - Authority"
$data.Range("D6").Value = "https://larimerco-ns5.trilynx-novastar.systems/novastar/operator"
$data.Range("E6").Value = "https://www.boxelderauthority.org/"
$data.Range("F6").Value = "Yes"
$data.Range("G6").Value = -105.08399
$data.Range("H6").Value = 40.40446

# ---------------------------------------------------------------------
# Notes sheet: replace the inventory description text. Done here so the
# shared-string table appends new unique strings in the same order as
# the authored edit.
# ---------------------------------------------------------------------
$notes.Range("A1").Value = "This Excel workbook is an initial inventory of stormwater/floodplain organizations that will be shown on the Current Conditions / Environment - Floods map."

$data.Range("K6").Value = "Coordinates are for consulting company."

# ---------------------------------------------------------------------
# Sheet view bookkeeping to match the saved selection state. The Notes
# sheet's stale D26 selection is cleared back to A1, and the Data sheet
# stays the active/frozen-pane sheet with the bottom-right pane's active
# cell moved to K6 (the new Note cell that was just filled in). Data is
# selected last so it remains the active tab.
# ---------------------------------------------------------------------
[void]$notes.Select()
[void]$notes.Range("A1").Select()

[void]$data.Select()
[void]$data.Range("K6").Select()
